$d = $word.ActiveDocument

# --- 1. Refresh the generated footer timestamp -----------------------------
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
$found = $footer.Range.Find.Execute("2025-06-30 12:12Z / ", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, "2025-07-02 02:48Z / ", 2)

# --- 2. Add the regression-test character styles: b, i, sub, sup, u --------
$wdStyleTypeCharacter = 2
$wdUnderlineSingle = 1

$newStyle = $d.Styles.Add("b", $wdStyleTypeCharacter)
$newStyle.BaseStyle = "DefaultParagraphFont"
$newStyle.Priority = 1
$newStyle.QuickStyle = $true
$newStyle.Font.Bold = $true

$newStyle = $d.Styles.Add("i", $wdStyleTypeCharacter)
$newStyle.BaseStyle = "DefaultParagraphFont"
$newStyle.Priority = 1
$newStyle.QuickStyle = $true
$newStyle.Font.Italic = $true

$newStyle = $d.Styles.Add("sub", $wdStyleTypeCharacter)
$newStyle.BaseStyle = "DefaultParagraphFont"
$newStyle.Priority = 1
$newStyle.QuickStyle = $true
$newStyle.Font.Subscript = $true

$newStyle = $d.Styles.Add("sup", $wdStyleTypeCharacter)
$newStyle.BaseStyle = "DefaultParagraphFont"
$newStyle.Priority = 1
$newStyle.QuickStyle = $true
$newStyle.Font.Superscript = $true

$newStyle = $d.Styles.Add("u", $wdStyleTypeCharacter)
$newStyle.BaseStyle = "DefaultParagraphFont"
$newStyle.Priority = 1
$newStyle.QuickStyle = $true
$newStyle.Font.Underline = $wdUnderlineSingle
